$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 7285
$ws.Range("I40").Value = 4998.3335
$ws.Range("K40").Value = 4998.3335
$ws.Range("M40").Value = -4823.3335

# Row 132
$ws.Range("H132").Value = 3457.7297
$ws.Range("I132").Value = 2767.3794
$ws.Range("K132").Value = 8302.138199999999
$ws.Range("M132").Value = -5772.138199999999

# Row 138
$ws.Range("H138").Value = 4363.273
$ws.Range("I138").Value = 3883.4285
$ws.Range("J138").Value = 4404.7407
$ws.Range("K138").Value = 11650.2855
$ws.Range("L138").Value = 13214.2221
$ws.Range("M138").Value = -6510.2855
$ws.Range("N138").Value = -23494.2221

# Row 141
$ws.Range("H141").Value = 7517.125
$ws.Range("I141").Value = 5763.8667
$ws.Range("K141").Value = 17291.6001
$ws.Range("M141").Value = -12111.6001

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1234.9333
$ws.Range("I45").Value = 1005.0909
$ws.Range("J45").Value = 1867
$ws.Range("K45").Value = 1005.0909
$ws.Range("L45").Value = 1867
$ws.Range("M45").Value = -628.0909
$ws.Range("N45").Value = -2621

# Row 63
$ws.Range("H63").Value = 3829
$ws.Range("I63").Value = 3829
$ws.Range("K63").Value = 3829
$ws.Range("M63").Value = -3143

# Row 66
$ws.Range("H66").Value = 3829
$ws.Range("I66").Value = 3829
$ws.Range("K66").Value = 19145
$ws.Range("M66").Value = -15713

# Row 74
$ws.Range("I74").Value = 2769.5
$ws.Range("J74").Value = 4218.4614
$ws.Range("K74").Value = 2769.5
$ws.Range("L74").Value = 4218.4614
$ws.Range("M74").Value = -1895.5
$ws.Range("N74").Value = -5966.4614

# Row 77
$ws.Range("I77").Value = 2769.5
$ws.Range("J77").Value = 4218.4614
$ws.Range("K77").Value = 13847.5
$ws.Range("L77").Value = 21092.307
$ws.Range("M77").Value = -9479.5
$ws.Range("N77").Value = -29828.307

# Row 88
$ws.Range("H88").Value = 1754.7778
$ws.Range("J88").Value = 2049.1428
$ws.Range("L88").Value = 2049.1428
$ws.Range("N88").Value = -2861.1428

# Row 91
$ws.Range("H91").Value = 1754.7778
$ws.Range("J91").Value = 2049.1428
$ws.Range("L91").Value = 2049.1428
$ws.Range("N91").Value = -4857.1428

# Row 125
$ws.Range("H125").Value = 199999
$ws.Range("J125").Value = 199999
$ws.Range("L125").Value = 199999
$ws.Range("N125").Value = -209839

# Row 132
$ws.Range("H132").Value = 12906.914
$ws.Range("J132").Value = 5876
$ws.Range("L132").Value = 17628
$ws.Range("N132").Value = -22688

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2530.6667
$ws.Range("I20").Value = 2429.84
$ws.Range("K20").Value = 2429.84
$ws.Range("M20").Value = -2182.84

# Row 99
$ws.Range("H99").Value = 3246.1
$ws.Range("I99").Value = 2682.75
$ws.Range("K99").Value = 2682.75
$ws.Range("M99").Value = -1184.75

# Row 105
$ws.Range("H105").Value = 2790.2856
$ws.Range("I105").Value = 2681.3635
$ws.Range("J105").Value = 3189.6667
$ws.Range("K105").Value = 2681.3635
$ws.Range("L105").Value = 3189.6667
$ws.Range("M105").Value = -934.3634999999999
$ws.Range("N105").Value = -6683.6667

# Row 134
$ws.Range("H134").Value = 3865.1035
$ws.Range("J134").Value = 7419.7144
$ws.Range("L134").Value = 22259.1432
$ws.Range("N134").Value = -27329.1432

$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Range("H28").Value = 10081
$ws.Range("J28").Value = 10081
$ws.Range("L28").Value = 10081
$ws.Range("N28").Value = -10571

# Row 31
$ws.Range("H31").Value = 3293.221
$ws.Range("I31").Value = 2675.7585
$ws.Range("K31").Value = 2675.7585
$ws.Range("M31").Value = -2380.7585

# Row 34
$ws.Range("H34").Value = 3293.221
$ws.Range("I34").Value = 2675.7585
$ws.Range("K34").Value = 2675.7585
$ws.Range("M34").Value = -2473.7585

# Row 50
$ws.Range("H50").Value = 49602
$ws.Range("I50").Value = 49999
$ws.Range("J50").Value = 49545.285
$ws.Range("K50").Value = 49999
$ws.Range("L50").Value = 49545.285
$ws.Range("M50").Value = -49374
$ws.Range("N50").Value = -50795.285

# Row 51
$ws.Range("H51").Value = 41255.75
$ws.Range("J51").Value = 39994
$ws.Range("L51").Value = 39994
$ws.Range("N51").Value = -41466

# Row 58
$ws.Range("H58").Value = 4203.6
$ws.Range("J58").Value = 9189.429
$ws.Range("L58").Value = 9189.429
$ws.Range("N58").Value = -9595.429

# Row 60
$ws.Range("H60").Value = 34531.535
$ws.Range("I60").Value = 11332.333
$ws.Range("J60").Value = 49997.668
$ws.Range("K60").Value = 11332.333
$ws.Range("L60").Value = 49997.668
$ws.Range("M60").Value = -10821.333
$ws.Range("N60").Value = -51019.668

# Row 61
$ws.Range("H61").Value = 41255.75
$ws.Range("J61").Value = 39994
$ws.Range("L61").Value = 39994
$ws.Range("N61").Value = -40690

# Row 62
$ws.Range("H62").Value = 10757.4
$ws.Range("I62").Value = 7812.4287
$ws.Range("J62").Value = 17629
$ws.Range("K62").Value = 7812.4287
$ws.Range("L62").Value = 17629
$ws.Range("M62").Value = -7188.4287
$ws.Range("N62").Value = -18877

# Row 65
$ws.Range("H65").Value = 10757.4
$ws.Range("I65").Value = 7812.4287
$ws.Range("J65").Value = 17629
$ws.Range("K65").Value = 39062.14350000001
$ws.Range("L65").Value = 88145
$ws.Range("M65").Value = -35942.14350000001
$ws.Range("N65").Value = -94385

# Row 136
$ws.Range("H136").Value = 4203.6
$ws.Range("J136").Value = 9189.429
$ws.Range("L136").Value = 27568.287
$ws.Range("N136").Value = -32668.287

$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 1239.8334
$ws.Range("I107").Value = 339.81818
$ws.Range("K107").Value = 1019.45454
$ws.Range("M107").Value = 900.54546

# Row 129
$ws.Range("H129").Value = 970.6111
$ws.Range("I129").Value = 807.625
$ws.Range("J129").Value = 1101
$ws.Range("K129").Value = 2422.875
$ws.Range("L129").Value = 3303
$ws.Range("M129").Value = 2577.125
$ws.Range("N129").Value = -13303

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 736579.4399999999
$ws.Range("I43").Value = 1255004.8
$ws.Range("J43").Value = 45345.668
$ws.Range("K43").Value = 1255004.8
$ws.Range("L43").Value = 45345.668
$ws.Range("M43").Value = -1254853.8
$ws.Range("N43").Value = -45647.668

# Row 70
$ws.Range("H70").Value = 5695.8
$ws.Range("I70").Value = 5633.3335
$ws.Range("J70").Value = 5789.5
$ws.Range("K70").Value = 5633.3335
$ws.Range("L70").Value = 5789.5
$ws.Range("M70").Value = -5363.3335
$ws.Range("N70").Value = -6329.5

# Row 73
$ws.Range("H73").Value = 5695.8
$ws.Range("I73").Value = 5633.3335
$ws.Range("J73").Value = 5789.5
$ws.Range("K73").Value = 5633.3335
$ws.Range("L73").Value = 5789.5
$ws.Range("M73").Value = -4697.3335
$ws.Range("N73").Value = -7661.5

# Row 126
$ws.Range("H126").Value = 103754.164
$ws.Range("I126").Value = 153001.12
$ws.Range("K126").Value = 459003.36
$ws.Range("M126").Value = -456533.36

# Row 132
$ws.Range("H132").Value = 2911.8772
$ws.Range("I132").Value = 2531.158
$ws.Range("J132").Value = 3673.3157
$ws.Range("K132").Value = 7593.474
$ws.Range("L132").Value = 11019.9471
$ws.Range("M132").Value = -5063.474
$ws.Range("N132").Value = -16079.9471

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 61544.047
$ws.Range("I40").Value = 76615.47
$ws.Range("K40").Value = 76615.47
$ws.Range("M40").Value = -76479.47

# Row 122
$ws.Range("H122").Value = 3476.6743
$ws.Range("J122").Value = 4216.9287
$ws.Range("L122").Value = 12650.7861
$ws.Range("N122").Value = -17550.7861

# Row 127
$ws.Range("H127").Value = 169999
$ws.Range("J127").Value = 169999
$ws.Range("L127").Value = 169999
$ws.Range("N127").Value = -179919

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1375
$ws.Range("I81").Value = 1500
$ws.Range("J81").Value = 1250
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 2500
$ws.Range("M81").Value = -1939
$ws.Range("N81").Value = -4622

# Row 84
$ws.Range("H84").Value = 1375
$ws.Range("I84").Value = 1500
$ws.Range("J84").Value = 1250
$ws.Range("K84").Value = 15000
$ws.Range("L84").Value = 12500
$ws.Range("M84").Value = -9696
$ws.Range("N84").Value = -23108

# Row 96
$ws.Range("H96").Value = 2892.923
$ws.Range("I96").Value = 2900.5
$ws.Range("K96").Value = 2900.5
$ws.Range("M96").Value = -1527.5

# Row 125
$ws.Range("H125").Value = 131666
$ws.Range("J125").Value = 131666
$ws.Range("L125").Value = 131666
$ws.Range("N125").Value = -141506

# Row 132
$ws.Range("H132").Value = 4222.116
$ws.Range("I132").Value = 3892.4075
$ws.Range("K132").Value = 11677.2225
$ws.Range("M132").Value = -9147.2225
